$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mso_tn_objects")

# --- Insert "sec_subnet" column (new I), shifting old gateway..domain_type right by one ---
$ws.Columns("I").Insert()
$ws.Range("I1").Value = "sec_subnet"
$ws.Range("I1").Style = "Normal"
$ws.Range("J1").Style = "Normal"

$ws.Range("I2").Value = "no"
$ws.Range("I3").Value = "no"
$ws.Range("I4").Value = "no"
$ws.Range("I5").Value = "no"
$ws.Range("I6").Value = "no"
$ws.Range("I7").Value = "no"

# --- Insert "gw_sec" column (new K), shifting old scope..domain_type right by one ---
$ws.Columns("K").Insert()
$ws.Range("K1").Value = "gw_sec"
$ws.Range("K1").Style = "Normal"

# --- Insert "host_route" column (new M), shifting old anp_name..domain_type right by one ---
$ws.Columns("M").Insert()
$ws.Range("M1").Value = "host_route"

$ws.Range("M2").Value = "no"
$ws.Range("M3").Value = "no"
$ws.Range("M4").Value = "no"
$ws.Range("M5").Value = "no"
$ws.Range("M6").Value = "no"
$ws.Range("M7").Value = "no"
$ws.Range("M8").Value = "no"
$ws.Range("M9").Value = "no"
$ws.Range("M8").Font.Color = $ws.Range("H8").Font.Color
$ws.Range("M9").Font.Color = $ws.Range("H9").Font.Color

# --- Match the bestFit column widths Excel would have computed for the new columns ---
$ws.Columns("I:I").ColumnWidth = 11.998697916666666
$ws.Columns("K:K").ColumnWidth = 8.998697916666666
$ws.Columns("M:M").ColumnWidth = 11.666666666666666

# --- Refresh the autofilter & used range over the new A1:Q9 extent ---
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:Q9").AutoFilter()

# --- Keep the named _FilterDatabase range for this sheet in sync with the new extent ---
foreach ($n in $wb.Names) {
  if ($n.Name -eq "mso_tn_objects!_FilterDatabase") {
    $n.RefersTo = "=mso_tn_objects!`$A`$1:`$Q`$9"
  }
}

# --- Activate the second sheet (mso_staticport), moving the tab-selected marker ---
$ws2 = $wb.Worksheets.Item("mso_staticport")
$ws2.Activate()
